$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 6).Value = 7957  # was 7930
$ws.Cells.Item(3, 6).Value = 110  # was 109
$ws.Cells.Item(4, 6).Value = 88  # was 86
$ws.Cells.Item(5, 6).Value = 19244  # was 17636
$ws.Cells.Item(8, 6).Value = 661  # was 657
$ws.Cells.Item(9, 6).Value = 455  # was 453
$ws.Cells.Item(10, 6).Value = 141  # was 140
$ws.Cells.Item(11, 6).Value = 443  # was 442
$ws.Cells.Item(12, 6).Value = 789  # was 787
$ws.Cells.Item(14, 6).Value = 591  # was 167
$ws.Cells.Item(15, 6).Value = 358  # was 354
$ws.Cells.Item(16, 6).Value = 23  # was 22
$ws.Cells.Item(18, 6).Value = 146  # was 145
$ws.Cells.Item(19, 6).Value = 402  # was 401
$ws.Cells.Item(20, 6).Value = 413  # was 412
$ws.Cells.Item(21, 6).Value = 1106  # was 1104
$ws.Cells.Item(22, 6).Value = 82  # was 81
$ws.Cells.Item(23, 6).Value = 663  # was 661
$ws.Cells.Item(24, 6).Value = 2262  # was 2247
$ws.Cells.Item(25, 6).Value = 781  # was 773
$ws.Cells.Item(29, 6).Value = 624  # was 623
$ws.Cells.Item(30, 6).Value = 563  # was 561

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(4, 6).Value = 333  # was 331

$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 6).Value = 496  # was 490

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 6).Value = 496  # was 490
$ws.Cells.Item(3, 6).Value = 7957  # was 7930
$ws.Cells.Item(4, 6).Value = 110  # was 109
$ws.Cells.Item(5, 6).Value = 88  # was 86
$ws.Cells.Item(7, 6).Value = 19246  # was 17642
$ws.Cells.Item(10, 6).Value = 661  # was 657
$ws.Cells.Item(11, 6).Value = 455  # was 453
$ws.Cells.Item(13, 6).Value = 141  # was 140
$ws.Cells.Item(14, 6).Value = 443  # was 442
$ws.Cells.Item(15, 6).Value = 333  # was 331
$ws.Cells.Item(18, 6).Value = 789  # was 787
$ws.Cells.Item(20, 6).Value = 591  # was 169
$ws.Cells.Item(21, 6).Value = 358  # was 354
$ws.Cells.Item(23, 6).Value = 23  # was 22
$ws.Cells.Item(28, 6).Value = 146  # was 145
$ws.Cells.Item(29, 6).Value = 402  # was 401
$ws.Cells.Item(30, 6).Value = 413  # was 412
$ws.Cells.Item(31, 6).Value = 1106  # was 1104
$ws.Cells.Item(32, 6).Value = 82  # was 81
$ws.Cells.Item(33, 6).Value = 663  # was 661
$ws.Cells.Item(34, 6).Value = 2262  # was 2248
$ws.Cells.Item(35, 6).Value = 781  # was 773
$ws.Cells.Item(40, 6).Value = 624  # was 623
$ws.Cells.Item(41, 6).Value = 563  # was 561
